$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 31   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/28/2024  Through  11/3/2024"

# --- Cells changing from numeric to text marker ("0" / "***.*") ---
$ws.Range("C31").Copy($ws.Range("C14"))
$ws.Range("C31").Copy($ws.Range("G15"))
$ws.Range("E31").Copy($ws.Range("H15"))
$ws.Range("C31").Copy($ws.Range("G27"))
$ws.Range("E31").Copy($ws.Range("H27"))

# --- Cells changing from text marker to numeric ---
$ws.Range("J31").Copy($ws.Range("D17"))
$ws.Range("D17").Value = 6
$ws.Range("K31").Copy($ws.Range("E17"))
$ws.Range("E17").Value = 33.333333333333
$ws.Range("J31").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 2
$ws.Range("K31").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("J31").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 1
$ws.Range("J31").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 2
$ws.Range("K31").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -50
$ws.Range("J31").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("K31").Copy($ws.Range("E29"))
$ws.Range("E29").Value = 300
$ws.Range("J31").Copy($ws.Range("G29"))
$ws.Range("G29").Value = 1
$ws.Range("K31").Copy($ws.Range("H29"))
$ws.Range("H29").Value = 400
$ws.Range("J31").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("K31").Copy($ws.Range("E30"))
$ws.Range("E30").Value = 100
$ws.Range("J31").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 1
$ws.Range("K31").Copy($ws.Range("H30"))
$ws.Range("H30").Value = 200

# --- Simple numeric value updates ---
$ws.Range("N15").Value = -76.315789473684
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -57.142857142857
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 143
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = -4.666666666666
$ws.Range("L16").Value = -4.026845637583
$ws.Range("M16").Value = -29.901960784313
$ws.Range("N16").Value = -80.544217687074
$ws.Range("C17").Value = 8
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 91.666666666666
$ws.Range("I17").Value = 257
$ws.Range("J17").Value = 210
$ws.Range("K17").Value = 22.380952380952
$ws.Range("L17").Value = 29.145728643216
$ws.Range("M17").Value = 62.658227848101
$ws.Range("N17").Value = -58.945686900958
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -33.333333333333
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 76
$ws.Range("J18").Value = 96
$ws.Range("K18").Value = -20.833333333333
$ws.Range("L18").Value = -49.333333333333
$ws.Range("M18").Value = -16.483516483516
$ws.Range("N18").Value = -89.645776566757
$ws.Range("C19").Value = 5
$ws.Range("E19").Value = -58.333333333333
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -8.108108108108
$ws.Range("I19").Value = 423
$ws.Range("J19").Value = 330
$ws.Range("K19").Value = 28.181818181818
$ws.Range("L19").Value = 21.902017291066
$ws.Range("M19").Value = 60.836501901140
$ws.Range("N19").Value = 4.444444444444
$ws.Range("D20").Value = 6
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -77.777777777777
$ws.Range("J20").Value = 63
$ws.Range("K20").Value = -22.222222222222
$ws.Range("N20").Value = -59.166666666666
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -47.058823529411
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = -3.75
$ws.Range("I21").Value = 962
$ws.Range("J21").Value = 859
$ws.Range("K21").Value = 11.990686845168
$ws.Range("L21").Value = 6.888888888888
$ws.Range("M21").Value = 28.437917222964
$ws.Range("N21").Value = -64.197990323781
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = -18.181818181818
$ws.Range("L22").Value = -40
$ws.Range("M22").Value = 63.636363636363
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = -71.428571428571
$ws.Range("J23").Value = 76
$ws.Range("K23").Value = -14.473684210526
$ws.Range("L23").Value = -16.666666666666
$ws.Range("M23").Value = 62.5
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 35.714285714285
$ws.Range("F24").Value = 124
$ws.Range("G24").Value = 85
$ws.Range("H24").Value = 45.882352941176
$ws.Range("I24").Value = 1089
$ws.Range("J24").Value = 1064
$ws.Range("K24").Value = 2.349624060150
$ws.Range("L24").Value = 11.577868852459
$ws.Range("M24").Value = 26.334106728538
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 70
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = 75
$ws.Range("I25").Value = 658
$ws.Range("J25").Value = 618
$ws.Range("K25").Value = 6.472491909385
$ws.Range("L25").Value = 18.558558558558
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 34
$ws.Range("H26").Value = 17.647058823529
$ws.Range("I26").Value = 383
$ws.Range("J26").Value = 349
$ws.Range("K26").Value = 9.742120343839
$ws.Range("L26").Value = 9.742120343839
$ws.Range("M26").Value = -12.756264236902
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 40
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = 33.333333333333
$ws.Range("L28").Value = -11.111111111111
$ws.Range("C29").Value = 4
$ws.Range("F29").Value = 5
$ws.Range("I29").Value = 19
$ws.Range("J29").Value = 10
$ws.Range("K29").Value = 90
$ws.Range("L29").Value = -20.833333333333
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -79.347826086956
$ws.Range("C30").Value = 2
$ws.Range("F30").Value = 3
$ws.Range("I30").Value = 14
$ws.Range("J30").Value = 9
$ws.Range("K30").Value = 55.555555555555
$ws.Range("L30").Value = -26.315789473684
$ws.Range("M30").Value = -17.647058823529
$ws.Range("N30").Value = -82.716049382716